$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 733.94116
$ws.Range("J17").Value = 733.94116
$ws.Range("L17").Value = 2201.82348
$ws.Range("N17").Value = -2537.82348

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 9552.666999999999
$ws.Range("I76").Value = 9372.267
$ws.Range("K76").Value = 9372.267
$ws.Range("M76").Value = -9057.267

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 9552.666999999999
$ws.Range("I79").Value = 9372.267
$ws.Range("K79").Value = 9372.267
$ws.Range("M79").Value = -8280.267

# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 113227
$ws.Range("I80").Value = 801
$ws.Range("K80").Value = 2403
$ws.Range("M80").Value = -1405

# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 113227
$ws.Range("I83").Value = 801
$ws.Range("K83").Value = 7209
$ws.Range("M83").Value = -2217

# Row 100: Asking for a Friend
$ws.Range("H100").Value = 3170.1428
$ws.Range("J100").Value = 2416.6667
$ws.Range("L100").Value = 2416.6667
$ws.Range("N100").Value = -3498.6667

# Row 118: Crafty Concoctions
$ws.Range("H118").Value = 591.2857
$ws.Range("I118").Value = 523.1667
$ws.Range("K118").Value = 1569.5001
$ws.Range("M118").Value = 87.49990000000003

# Row 135: For Tired Minds
$ws.Range("H135").Value = 5818.2173
$ws.Range("I135").Value = 1097.7368
$ws.Range("K135").Value = 9879.6312
$ws.Range("M135").Value = -7344.6312

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2309.8
$ws.Range("I137").Value = 1120.5
$ws.Range("K137").Value = 3361.5
$ws.Range("M137").Value = -811.5

$ws = $wb.Worksheets.Item("ARM")
# Row 11: Rodents of Unusual Size
$ws.Range("H11").Value = 33336334
$ws.Range("I11").Value = 25004502
$ws.Range("K11").Value = 25004502
$ws.Range("M11").Value = -25004358

# Row 32: Ingot We Trust
$ws.Range("H32").Value = 2003.3582
$ws.Range("I32").Value = 800.4219000000001
$ws.Range("K32").Value = 800.4219000000001
$ws.Range("M32").Value = -513.4219000000001

# Row 88: The Mast Chance
$ws.Range("H88").Value = 1311.5
$ws.Range("I88").Value = 787.25
$ws.Range("J88").Value = 1661
$ws.Range("K88").Value = 787.25
$ws.Range("L88").Value = 1661
$ws.Range("M88").Value = -381.25
$ws.Range("N88").Value = -2473

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 1311.5
$ws.Range("I91").Value = 787.25
$ws.Range("J91").Value = 1661
$ws.Range("K91").Value = 787.25
$ws.Range("L91").Value = 1661
$ws.Range("M91").Value = 616.75
$ws.Range("N91").Value = -4469

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 6039.5
$ws.Range("I102").Value = 5709.3335
$ws.Range("K102").Value = 5709.3335
$ws.Range("M102").Value = -4087.3335

# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 2674.6667
$ws.Range("I110").Value = 2050.4443
$ws.Range("J110").Value = 3611
$ws.Range("K110").Value = 2050.4443
$ws.Range("L110").Value = 3611
$ws.Range("M110").Value = -5.444300000000112
$ws.Range("N110").Value = -7701

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Unbreaker
$ws.Range("H80").Value = 2775.2222
$ws.Range("I80").Value = 2664.6667
$ws.Range("K80").Value = 2664.6667
$ws.Range("M80").Value = -1666.6667

# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 2775.2222
$ws.Range("I83").Value = 2664.6667
$ws.Range("K83").Value = 13323.3335
$ws.Range("M83").Value = -8331.333500000001

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 3120.524
$ws.Range("I99").Value = 2913.5881
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2913.5881
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -1415.5881
$ws.Range("N99").Value = -6996

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 5750
$ws.Range("I105").Value = 5750
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5750
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -4003

# Row 107: The Gold Experience
$ws.Range("H107").Value = 1989.2433
$ws.Range("I107").Value = 1847.4117
$ws.Range("J107").Value = 3596.6667
$ws.Range("K107").Value = 1847.4117
$ws.Range("L107").Value = 3596.6667
$ws.Range("M107").Value = 72.58829999999989
$ws.Range("N107").Value = -7436.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1951.125
$ws.Range("I105").Value = 1918.1666
$ws.Range("J105").Value = 2050
$ws.Range("K105").Value = 1918.1666
$ws.Range("L105").Value = 2050
$ws.Range("M105").Value = -171.1666
$ws.Range("N105").Value = -5544

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 66598.5
$ws.Range("J141").Value = 66598.5
$ws.Range("L141").Value = 66598.5
$ws.Range("N141").Value = -76958.5

$ws = $wb.Worksheets.Item("CUL")
# Row 86: Let's Not Get Sappy
$ws.Range("H86").Value = 8705.77
$ws.Range("J86").Value = 17498
$ws.Range("L86").Value = 52494
$ws.Range("N86").Value = -54866

# Row 89: Luxury Spillover (L)
$ws.Range("H89").Value = 8705.77
$ws.Range("J89").Value = 17498
$ws.Range("L89").Value = 157482
$ws.Range("N89").Value = -169338

# Row 112: Sweet Tooth
$ws.Range("H112").Value = 14712311
$ws.Range("J112").Value = 16673500
$ws.Range("L112").Value = 50020500
$ws.Range("N112").Value = -50022716

# Row 115: Mixology
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

# Row 124: Bobbing for Compliments
$ws.Range("H124").Value = 7715
$ws.Range("I124").Value = 7715
$ws.Range("K124").Value = 23145
$ws.Range("M124").Value = -18235

# Row 126: Imperial Palate
$ws.Range("H126").Value = 8129.8
$ws.Range("I126").Value = 2662.25
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 7986.75
$ws.Range("L126").Value = 90000
$ws.Range("M126").Value = -3046.75
$ws.Range("N126").Value = -99880

$ws = $wb.Worksheets.Item("LTW")
# Row 20: Choke Hold
$ws.Range("H20").Value = 8335000
$ws.Range("J20").Value = 5000
$ws.Range("L20").Value = 5000
$ws.Range("N20").Value = -5452

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 4244.857
$ws.Range("I46").Value = 4431.5
$ws.Range("K46").Value = 4431.5
$ws.Range("M46").Value = -4243.5

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 950.4286
$ws.Range("I93").Value = 950.4286
$ws.Range("K93").Value = 950.4286
$ws.Range("M93").Value = 297.5714

$ws = $wb.Worksheets.Item("WVR")
# Row 74: Clothing the Naked Truth
$ws.Range("H74").Value = 20974.5
$ws.Range("J74").Value = 20974.5
$ws.Range("L74").Value = 20974.5
$ws.Range("N74").Value = -22846.5

# Row 77: When in Robes (L)
$ws.Range("H77").Value = 20974.5
$ws.Range("J77").Value = 20974.5
$ws.Range("L77").Value = 62923.5
$ws.Range("N77").Value = -72283.5

# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 3356.7222
$ws.Range("I81").Value = 3497.818
$ws.Range("J81").Value = 3135
$ws.Range("K81").Value = 6995.636
$ws.Range("L81").Value = 6270
$ws.Range("M81").Value = -5934.636
$ws.Range("N81").Value = -8392

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 3356.7222
$ws.Range("I84").Value = 3497.818
$ws.Range("J84").Value = 3135
$ws.Range("K84").Value = 34978.18
$ws.Range("L84").Value = 31350
$ws.Range("M84").Value = -29674.18
$ws.Range("N84").Value = -41958

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 2795.4722
$ws.Range("I132").Value = 3370.0908
$ws.Range("J132").Value = 1892.5
$ws.Range("K132").Value = 10110.2724
$ws.Range("L132").Value = 5677.5
$ws.Range("M132").Value = -7580.2724
$ws.Range("N132").Value = -10737.5
